$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column E (Population) before the existing Latitude column.
# Using EntireColumn.Insert shifts the old E:L data (and their per-cell
# styles / column formatting) over to F:M, matching Excel's native
# "Insert Column" behaviour.
$ws.Columns("E").Insert()

# Give the new column the same (un-autofit) width as its neighbours.
$ws.Columns("E").ColumnWidth = 10.33

# Header for the new column
$ws.Range("E1").Value = "Population"

# Population figures (millions) for each PoP row
$ws.Range("E2").Value = 61
$ws.Range("E3").Value = 81
$ws.Range("E4").Value = 11
$ws.Range("E5").Value = 30
$ws.Range("E6").Value = 10
$ws.Range("E7").Value = 6
$ws.Range("E8").Value = 5
$ws.Range("E9").Value = 23
$ws.Range("E10").Value = 5
$ws.Range("E11").Value = 10
$ws.Range("E12").Value = 10
$ws.Range("E13").Value = 9
$ws.Range("E14").Value = 23
$ws.Range("E15").Value = 8

# Fix the Country/City values that were swapped for the ICCS row.
$ws.Range("C10").Value = "Greece"
$ws.Range("D10").Value = "Athens"

# New broker population row for NORDUNET.
$ws.Range("C16").Value = "Denmark"
$ws.Range("E16").Value = 6

# Match the author's final selection.
$ws.Range("F19").Select()
